$d = $word.ActiveDocument

# --- Title paragraph: collapse the split "WS"/"2022/"/"2023" runs (no text change) ---
$d.Content.Find.Execute("Platforms and Systems for eLearning WS 2022/2023 ", $true, $false, $false, $false, $false, $true, 1, $false, "Platforms and Systems for eLearning WS 2022/2023 ", 2)

# --- "Project 3" paragraph: collapse the split runs (no text change) ---
$d.Content.Find.Execute("Project 3", $true, $false, $false, $false, $false, $true, 1, $false, "Project 3", 2)

# --- Andreas Scholl paragraph: the trailing tab becomes real task text ---
$d.Content.Find.Execute("Andreas Scholl: " + [char]9, $true, $false, $false, $false, $false, $true, 1, $false, "Andreas Scholl: Basic structure (tabs) and plots", 2)

# --- Katrin Peikert paragraph: collapse the split runs (no text change) ---
$d.Content.Find.Execute("Katrin Peikert: Wordcloud", $true, $false, $false, $false, $false, $true, 1, $false, "Katrin Peikert: Wordcloud", 2)

# --- Kenny La paragraph: add trailing space + new task description ---
$d.Content.Find.Execute("Kenny La:", $true, $false, $false, $false, $false, $true, 1, $false, "Kenny La: Mapclient, Citiymap, Heatmap", 2)
